$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update header text: June 16 -> June 17
$ws.Name = "Through 2022-06-17"
$ws.Range("B1").Value = "June 2022 (through June 17)"

# Update/add cell counts (carjacking data additions for 2022-06-25)
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 7
$ws.Range("T3").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("Z4").Value = 6
$ws.Range("AR4").Value = 1
$ws.Range("T5").Value = 2
$ws.Range("B9").Value = 4
$ws.Range("AL9").Value = 2
$ws.Range("Z10").Value = 2
$ws.Range("AF10").Value = 1
$ws.Range("B12").Value = 4
$ws.Range("B14").Value = 6
$ws.Range("H14").Value = 7
$ws.Range("B15").Value = 2
$ws.Range("H15").Value = 3
$ws.Range("AF19").Value = 1
$ws.Range("T23").Value = 2
$ws.Range("H33").Value = 2
$ws.Range("Z36").Value = 1
$ws.Range("N39").Value = 1
$ws.Range("AF41").Value = 1
$ws.Range("H47").Value = 1
$ws.Range("N54").Value = 1
$ws.Range("B68").Value = 2
$ws.Range("H68").Value = 3
$ws.Range("N75").Value = 2
$ws.Range("H77").Value = 1
$ws.Range("H86").Value = 1
$ws.Range("B92").Value = 3
